# Update Name of Algo
# Applies the numeric corrections to the RandomForest imputation results
# on Sheet1, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = -7.321699999999993
$ws.Range("A12").Value = -21.4147
$ws.Range("D23").Value = -8.2957
$ws.Range("D28").Value = -8.272800000000002
$ws.Range("A32").Value = -21.4918
$ws.Range("D32").Value = -6.733499999999993
$ws.Range("D34").Value = -8.128399999999997
$ws.Range("A36").Value = -20.0503
$ws.Range("A38").Value = -19.49599999999998
$ws.Range("D42").Value = -9.145999999999994
$ws.Range("A46").Value = -21.9308
$ws.Range("A54").Value = -22.129
$ws.Range("D54").Value = -7.847900000000001
$ws.Range("A55").Value = -22.1489
$ws.Range("A67").Value = -21.44199999999996
$ws.Range("A69").Value = -21.52639999999997
$ws.Range("A72").Value = -21.84949999999999
$ws.Range("A91").Value = -20.36239999999998
$ws.Range("D97").Value = -8.390699999999995
$ws.Range("A99").Value = -21.8191
$ws.Range("D99").Value = -7.828200000000001
$ws.Range("D101").Value = -8.121199999999989
$ws.Range("A104").Value = -21.48539999999998
